$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 8).Value = "kitchens"
$ws.Cells.Item(2, 9).Value = "target"
$ws.Cells.Item(2, 11).Value = "j"
$ws.Cells.Item(2, 12).Value = "stimuli/img_k3abb.png"
$ws.Cells.Item(2, 13).Value = 35.54054054054054
$ws.Cells.Item(2, 14).Value = 16.54054054054054
$ws.Cells.Item(2, 15).Value = 26.04054054054054
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = 1
$ws.Cells.Item(2, 19).Value = 1
$ws.Cells.Item(2, 20).Value = 1
$ws.Cells.Item(2, 21).Value = 1
$ws.Cells.Item(2, 22).Value = 1
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 8).Value = "bedrooms"
$ws.Cells.Item(3, 9).Value = "distractor"
$ws.Cells.Item(3, 11).Value = "f"
$ws.Cells.Item(3, 12).Value = "stimuli/img_088sa.png"
$ws.Cells.Item(3, 13).Value = 67.36363636363636
$ws.Cells.Item(3, 14).Value = 42.12121212121212
$ws.Cells.Item(3, 15).Value = 54.74242424242424
$ws.Cells.Item(3, 16).Value = 33
$ws.Cells.Item(3, 17).Value = 4
$ws.Cells.Item(3, 18).Value = 4
$ws.Cells.Item(3, 19).Value = 4
$ws.Cells.Item(3, 20).Value = 4
$ws.Cells.Item(3, 21).Value = 4
$ws.Cells.Item(3, 22).Value = 4
$ws.Cells.Item(4, 3).Value = 2
$ws.Cells.Item(4, 8).Value = "kitchens"
$ws.Cells.Item(4, 9).Value = "target"
$ws.Cells.Item(4, 11).Value = "j"
$ws.Cells.Item(4, 12).Value = "stimuli/img_qbkdt.png"
$ws.Cells.Item(4, 13).Value = 69.45714285714286
$ws.Cells.Item(4, 14).Value = 50.02857142857143
$ws.Cells.Item(4, 15).Value = 59.74285714285715
$ws.Cells.Item(4, 16).Value = 35
$ws.Cells.Item(4, 17).Value = 5
$ws.Cells.Item(4, 18).Value = 5
$ws.Cells.Item(4, 19).Value = 5
$ws.Cells.Item(4, 20).Value = 5
$ws.Cells.Item(4, 21).Value = 5
$ws.Cells.Item(4, 22).Value = 5
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 8).Value = "kitchens"
$ws.Cells.Item(5, 9).Value = "target"
$ws.Cells.Item(5, 11).Value = "j"
$ws.Cells.Item(5, 12).Value = "stimuli/img_faly8.png"
$ws.Cells.Item(5, 13).Value = 33.41176470588236
$ws.Cells.Item(5, 14).Value = 19.23529411764706
$ws.Cells.Item(5, 15).Value = 26.32352941176471
$ws.Cells.Item(5, 16).Value = 34
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = 1
$ws.Cells.Item(5, 19).Value = 1
$ws.Cells.Item(5, 20).Value = 1
$ws.Cells.Item(5, 21).Value = 1
$ws.Cells.Item(5, 22).Value = 1
$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 8).Value = "bedrooms"
$ws.Cells.Item(6, 9).Value = "distractor"
$ws.Cells.Item(6, 11).Value = "f"
$ws.Cells.Item(6, 12).Value = "stimuli/img_71mhq.png"
$ws.Cells.Item(6, 13).Value = 69.34210526315789
$ws.Cells.Item(6, 14).Value = 47.02631578947368
$ws.Cells.Item(6, 15).Value = 58.18421052631579
$ws.Cells.Item(6, 16).Value = 38
$ws.Cells.Item(6, 17).Value = 5
$ws.Cells.Item(6, 18).Value = 5
$ws.Cells.Item(6, 19).Value = 5
$ws.Cells.Item(6, 20).Value = 5
$ws.Cells.Item(6, 21).Value = 5
$ws.Cells.Item(6, 22).Value = 5
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 8).Value = "living_rooms"
$ws.Cells.Item(7, 12).Value = "stimuli/img_4974k.png"
$ws.Cells.Item(7, 13).Value = 79.89130434782609
$ws.Cells.Item(7, 14).Value = 58.19565217391305
$ws.Cells.Item(7, 15).Value = 69.04347826086956
$ws.Cells.Item(7, 16).Value = 46
$ws.Cells.Item(8, 3).Value = 2
$ws.Cells.Item(8, 8).Value = "kitchens"
$ws.Cells.Item(8, 9).Value = "target"
$ws.Cells.Item(8, 11).Value = "j"
$ws.Cells.Item(8, 12).Value = "stimuli/img_hfz8w.png"
$ws.Cells.Item(8, 13).Value = 55.46153846153846
$ws.Cells.Item(8, 14).Value = 27.28205128205128
$ws.Cells.Item(8, 15).Value = 41.37179487179487
$ws.Cells.Item(8, 16).Value = 39
$ws.Cells.Item(8, 17).Value = 2
$ws.Cells.Item(8, 18).Value = 2
$ws.Cells.Item(8, 19).Value = 2
$ws.Cells.Item(8, 20).Value = 2
$ws.Cells.Item(8, 21).Value = 2
$ws.Cells.Item(8, 22).Value = 2
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 8).Value = "kitchens"
$ws.Cells.Item(9, 9).Value = "target"
$ws.Cells.Item(9, 11).Value = "j"
$ws.Cells.Item(9, 12).Value = "stimuli/img_5949k.png"
$ws.Cells.Item(9, 13).Value = 60.8
$ws.Cells.Item(9, 14).Value = 39.2
$ws.Cells.Item(9, 15).Value = 50
$ws.Cells.Item(9, 16).Value = 35
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(9, 18).Value = 3
$ws.Cells.Item(9, 19).Value = 3
$ws.Cells.Item(9, 20).Value = 3
$ws.Cells.Item(9, 21).Value = 3
$ws.Cells.Item(9, 22).Value = 3
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 8).Value = "bedrooms"
$ws.Cells.Item(10, 12).Value = "stimuli/img_d3t0o.png"
$ws.Cells.Item(10, 13).Value = 66.95121951219512
$ws.Cells.Item(10, 14).Value = 42.92682926829269
$ws.Cells.Item(10, 15).Value = 54.9390243902439
$ws.Cells.Item(10, 16).Value = 41
$ws.Cells.Item(10, 17).Value = 4
$ws.Cells.Item(10, 18).Value = 4
$ws.Cells.Item(10, 19).Value = 4
$ws.Cells.Item(10, 20).Value = 4
$ws.Cells.Item(10, 21).Value = 4
$ws.Cells.Item(10, 22).Value = 4
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 8).Value = "bedrooms"
$ws.Cells.Item(11, 9).Value = "distractor"
$ws.Cells.Item(11, 11).Value = "f"
$ws.Cells.Item(11, 12).Value = "stimuli/img_mqnl6.png"
$ws.Cells.Item(11, 13).Value = 70.7560975609756
$ws.Cells.Item(11, 14).Value = 45.68292682926829
$ws.Cells.Item(11, 15).Value = 58.21951219512195
$ws.Cells.Item(11, 16).Value = 41
$ws.Cells.Item(11, 17).Value = 5
$ws.Cells.Item(11, 18).Value = 5
$ws.Cells.Item(11, 19).Value = 5
$ws.Cells.Item(11, 20).Value = 5
$ws.Cells.Item(11, 21).Value = 5
$ws.Cells.Item(11, 22).Value = 5
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 8).Value = "living_rooms"
$ws.Cells.Item(12, 9).Value = "distractor"
$ws.Cells.Item(12, 11).Value = "f"
$ws.Cells.Item(12, 12).Value = "stimuli/img_lpas9.png"
$ws.Cells.Item(12, 13).Value = 59.36585365853659
$ws.Cells.Item(12, 14).Value = 39.09756097560975
$ws.Cells.Item(12, 15).Value = 49.23170731707317
$ws.Cells.Item(12, 16).Value = 41
$ws.Cells.Item(12, 17).Value = 4
$ws.Cells.Item(12, 18).Value = 4
$ws.Cells.Item(12, 19).Value = 4
$ws.Cells.Item(12, 20).Value = 4
$ws.Cells.Item(12, 21).Value = 4
$ws.Cells.Item(12, 22).Value = 4
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 8).Value = "bedrooms"
$ws.Cells.Item(13, 12).Value = "stimuli/img_bdz92.png"
$ws.Cells.Item(13, 13).Value = 63.72222222222222
$ws.Cells.Item(13, 14).Value = 42.63888888888889
$ws.Cells.Item(13, 15).Value = 53.18055555555556
$ws.Cells.Item(13, 16).Value = 36
$ws.Cells.Item(13, 17).Value = 4
$ws.Cells.Item(13, 18).Value = 4
$ws.Cells.Item(13, 19).Value = 4
$ws.Cells.Item(13, 20).Value = 4
$ws.Cells.Item(13, 21).Value = 4
$ws.Cells.Item(13, 22).Value = 4
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 12).Value = "stimuli/img_8fpog.png"
$ws.Cells.Item(14, 13).Value = 85.41666666666667
$ws.Cells.Item(14, 14).Value = 72.30555555555556
$ws.Cells.Item(14, 15).Value = 78.86111111111111
$ws.Cells.Item(14, 16).Value = 36
$ws.Cells.Item(14, 17).Value = 10
$ws.Cells.Item(14, 18).Value = 10
$ws.Cells.Item(14, 19).Value = 10
$ws.Cells.Item(14, 20).Value = 10
$ws.Cells.Item(14, 21).Value = 10
$ws.Cells.Item(14, 22).Value = 10
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 12).Value = "stimuli/img_l9t30.png"
$ws.Cells.Item(15, 13).Value = 67.2
$ws.Cells.Item(15, 14).Value = 43.14285714285715
$ws.Cells.Item(15, 15).Value = 55.17142857142858
$ws.Cells.Item(15, 16).Value = 35
$ws.Cells.Item(15, 17).Value = 4
$ws.Cells.Item(15, 18).Value = 4
$ws.Cells.Item(15, 19).Value = 4
$ws.Cells.Item(15, 20).Value = 4
$ws.Cells.Item(15, 21).Value = 4
$ws.Cells.Item(15, 22).Value = 4
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 8).Value = "bedrooms"
$ws.Cells.Item(16, 9).Value = "distractor"
$ws.Cells.Item(16, 11).Value = "f"
$ws.Cells.Item(16, 12).Value = "stimuli/img_uttnz.png"
$ws.Cells.Item(16, 13).Value = 69.91891891891892
$ws.Cells.Item(16, 14).Value = 49.91891891891892
$ws.Cells.Item(16, 15).Value = 59.91891891891892
$ws.Cells.Item(16, 16).Value = 37
$ws.Cells.Item(16, 17).Value = 5
$ws.Cells.Item(16, 18).Value = 5
$ws.Cells.Item(16, 19).Value = 5
$ws.Cells.Item(16, 20).Value = 5
$ws.Cells.Item(16, 21).Value = 5
$ws.Cells.Item(16, 22).Value = 5
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 8).Value = "kitchens"
$ws.Cells.Item(17, 9).Value = "target"
$ws.Cells.Item(17, 11).Value = "j"
$ws.Cells.Item(17, 12).Value = "stimuli/img_u9f9l.png"
$ws.Cells.Item(17, 13).Value = 77.78571428571429
$ws.Cells.Item(17, 14).Value = 57.25
$ws.Cells.Item(17, 15).Value = 67.51785714285714
$ws.Cells.Item(17, 16).Value = 28
$ws.Cells.Item(17, 17).Value = 7
$ws.Cells.Item(17, 18).Value = 7
$ws.Cells.Item(17, 19).Value = 7
$ws.Cells.Item(17, 20).Value = 7
$ws.Cells.Item(17, 21).Value = 7
$ws.Cells.Item(17, 22).Value = 7
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 12).Value = "stimuli/img_qmand.png"
$ws.Cells.Item(18, 13).Value = 86.11764705882354
$ws.Cells.Item(18, 14).Value = 71.02941176470588
$ws.Cells.Item(18, 15).Value = 78.57352941176471
$ws.Cells.Item(18, 16).Value = 34
$ws.Cells.Item(18, 17).Value = 10
$ws.Cells.Item(18, 18).Value = 10
$ws.Cells.Item(18, 19).Value = 10
$ws.Cells.Item(18, 20).Value = 10
$ws.Cells.Item(18, 21).Value = 10
$ws.Cells.Item(18, 22).Value = 10
$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(19, 8).Value = "living_rooms"
$ws.Cells.Item(19, 9).Value = "distractor"
$ws.Cells.Item(19, 11).Value = "f"
$ws.Cells.Item(19, 12).Value = "stimuli/img_xu1p3.png"
$ws.Cells.Item(19, 13).Value = 75.27659574468085
$ws.Cells.Item(19, 14).Value = 56.68085106382978
$ws.Cells.Item(19, 15).Value = 65.97872340425532
$ws.Cells.Item(19, 16).Value = 47
$ws.Cells.Item(19, 17).Value = 7
$ws.Cells.Item(19, 18).Value = 7
$ws.Cells.Item(19, 19).Value = 7
$ws.Cells.Item(19, 20).Value = 6
$ws.Cells.Item(19, 21).Value = 6
$ws.Cells.Item(19, 22).Value = 6
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 8).Value = "bedrooms"
$ws.Cells.Item(20, 9).Value = "distractor"
$ws.Cells.Item(20, 11).Value = "f"
$ws.Cells.Item(20, 12).Value = "stimuli/img_twj5p.png"
$ws.Cells.Item(20, 13).Value = 67.71739130434783
$ws.Cells.Item(20, 14).Value = 42.08695652173913
$ws.Cells.Item(20, 15).Value = 54.90217391304348
$ws.Cells.Item(20, 16).Value = 46
$ws.Cells.Item(20, 17).Value = 4
$ws.Cells.Item(20, 18).Value = 4
$ws.Cells.Item(20, 19).Value = 4
$ws.Cells.Item(20, 20).Value = 4
$ws.Cells.Item(20, 22).Value = 4
$ws.Cells.Item(21, 3).Value = 2
$ws.Cells.Item(21, 8).Value = "bedrooms"
$ws.Cells.Item(21, 9).Value = "distractor"
$ws.Cells.Item(21, 11).Value = "f"
$ws.Cells.Item(21, 12).Value = "stimuli/img_h13c3.png"
$ws.Cells.Item(21, 13).Value = 71.80555555555556
$ws.Cells.Item(21, 14).Value = 47.86111111111111
$ws.Cells.Item(21, 15).Value = 59.83333333333334
$ws.Cells.Item(21, 16).Value = 36
$ws.Cells.Item(21, 17).Value = 5
$ws.Cells.Item(21, 18).Value = 5
$ws.Cells.Item(21, 19).Value = 5
$ws.Cells.Item(21, 20).Value = 5
$ws.Cells.Item(21, 21).Value = 5
$ws.Cells.Item(21, 22).Value = 5
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 8).Value = "kitchens"
$ws.Cells.Item(22, 9).Value = "target"
$ws.Cells.Item(22, 11).Value = "j"
$ws.Cells.Item(22, 12).Value = "stimuli/img_01w8b.png"
$ws.Cells.Item(22, 13).Value = 78.91891891891892
$ws.Cells.Item(22, 14).Value = 61.21621621621622
$ws.Cells.Item(22, 15).Value = 70.06756756756756
$ws.Cells.Item(22, 16).Value = 37
$ws.Cells.Item(22, 17).Value = 8
$ws.Cells.Item(22, 18).Value = 8
$ws.Cells.Item(22, 19).Value = 8
$ws.Cells.Item(22, 20).Value = 8
$ws.Cells.Item(22, 21).Value = 8
$ws.Cells.Item(22, 22).Value = 8
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 12).Value = "stimuli/img_b89t4.png"
$ws.Cells.Item(23, 13).Value = 71.41463414634147
$ws.Cells.Item(23, 14).Value = 47.85365853658536
$ws.Cells.Item(23, 15).Value = 59.63414634146342
$ws.Cells.Item(23, 16).Value = 41
$ws.Cells.Item(23, 17).Value = 5
$ws.Cells.Item(23, 18).Value = 5
$ws.Cells.Item(23, 19).Value = 5
$ws.Cells.Item(23, 20).Value = 5
$ws.Cells.Item(23, 21).Value = 5
$ws.Cells.Item(23, 22).Value = 5
$ws.Cells.Item(24, 3).Value = 2
$ws.Cells.Item(24, 12).Value = "stimuli/img_0j24m.png"
$ws.Cells.Item(24, 13).Value = 63.6969696969697
$ws.Cells.Item(24, 14).Value = 35.75757575757576
$ws.Cells.Item(24, 15).Value = 49.72727272727273
$ws.Cells.Item(24, 16).Value = 33
$ws.Cells.Item(24, 17).Value = 3
$ws.Cells.Item(24, 18).Value = 3
$ws.Cells.Item(24, 19).Value = 3
$ws.Cells.Item(24, 20).Value = 3
$ws.Cells.Item(24, 21).Value = 3
$ws.Cells.Item(24, 22).Value = 3
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 8).Value = "bedrooms"
$ws.Cells.Item(25, 12).Value = "stimuli/img_huisn.png"
$ws.Cells.Item(25, 13).Value = 73.63888888888889
$ws.Cells.Item(25, 14).Value = 46.36111111111111
$ws.Cells.Item(25, 15).Value = 60
$ws.Cells.Item(25, 16).Value = 36
$ws.Cells.Item(25, 17).Value = 5
$ws.Cells.Item(25, 18).Value = 5
$ws.Cells.Item(25, 19).Value = 5
$ws.Cells.Item(25, 20).Value = 5
$ws.Cells.Item(25, 21).Value = 5
$ws.Cells.Item(25, 22).Value = 5
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 8).Value = "kitchens"
$ws.Cells.Item(26, 9).Value = "target"
$ws.Cells.Item(26, 11).Value = "j"
$ws.Cells.Item(26, 12).Value = "stimuli/img_as3da.png"
$ws.Cells.Item(26, 13).Value = 84.53125
$ws.Cells.Item(26, 14).Value = 63
$ws.Cells.Item(26, 15).Value = 73.765625
$ws.Cells.Item(26, 16).Value = 32
$ws.Cells.Item(26, 17).Value = 9
$ws.Cells.Item(26, 18).Value = 9
$ws.Cells.Item(26, 19).Value = 9
$ws.Cells.Item(26, 20).Value = 9
$ws.Cells.Item(26, 21).Value = 9
$ws.Cells.Item(26, 22).Value = 9
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 12).Value = "stimuli/img_ensho.png"
$ws.Cells.Item(27, 13).Value = 72.7948717948718
$ws.Cells.Item(27, 14).Value = 54.56410256410256
$ws.Cells.Item(27, 15).Value = 63.67948717948718
$ws.Cells.Item(27, 16).Value = 39
$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(28, 8).Value = "kitchens"
$ws.Cells.Item(28, 9).Value = "target"
$ws.Cells.Item(28, 11).Value = "j"
$ws.Cells.Item(28, 12).Value = "stimuli/img_fhm45.png"
$ws.Cells.Item(28, 13).Value = 76.75
$ws.Cells.Item(28, 14).Value = 57.71875
$ws.Cells.Item(28, 15).Value = 67.234375
$ws.Cells.Item(28, 16).Value = 32
$ws.Cells.Item(28, 17).Value = 7
$ws.Cells.Item(28, 18).Value = 7
$ws.Cells.Item(28, 19).Value = 7
$ws.Cells.Item(28, 20).Value = 7
$ws.Cells.Item(28, 21).Value = 7
$ws.Cells.Item(28, 22).Value = 7
$ws.Cells.Item(29, 3).Value = 2
$ws.Cells.Item(29, 12).Value = "stimuli/img_uspja.png"
$ws.Cells.Item(29, 13).Value = 54.90909090909091
$ws.Cells.Item(29, 14).Value = 29.12121212121212
$ws.Cells.Item(29, 15).Value = 42.01515151515152
$ws.Cells.Item(29, 16).Value = 33
$ws.Cells.Item(29, 17).Value = 2
$ws.Cells.Item(29, 18).Value = 2
$ws.Cells.Item(29, 19).Value = 2
$ws.Cells.Item(29, 20).Value = 2
$ws.Cells.Item(29, 21).Value = 2
$ws.Cells.Item(29, 22).Value = 2
$ws.Cells.Item(30, 3).Value = 2
$ws.Cells.Item(30, 8).Value = "living_rooms"
$ws.Cells.Item(30, 9).Value = "distractor"
$ws.Cells.Item(30, 11).Value = "f"
$ws.Cells.Item(30, 12).Value = "stimuli/img_g9od8.png"
$ws.Cells.Item(30, 13).Value = 59.34883720930232
$ws.Cells.Item(30, 14).Value = 37.83720930232558
$ws.Cells.Item(30, 15).Value = 48.59302325581395
$ws.Cells.Item(30, 16).Value = 43
$ws.Cells.Item(30, 17).Value = 4
$ws.Cells.Item(30, 18).Value = 4
$ws.Cells.Item(30, 19).Value = 4
$ws.Cells.Item(30, 20).Value = 4
$ws.Cells.Item(30, 21).Value = 4
$ws.Cells.Item(30, 22).Value = 4
$ws.Cells.Item(31, 3).Value = 2
$ws.Cells.Item(31, 12).Value = "stimuli/img_02alv.png"
$ws.Cells.Item(31, 13).Value = 61.8
$ws.Cells.Item(31, 14).Value = 37.8
$ws.Cells.Item(31, 15).Value = 49.8
$ws.Cells.Item(31, 16).Value = 45
$ws.Cells.Item(31, 17).Value = 4
$ws.Cells.Item(31, 18).Value = 4
$ws.Cells.Item(31, 19).Value = 4
$ws.Cells.Item(31, 20).Value = 4
$ws.Cells.Item(31, 21).Value = 5
$ws.Cells.Item(31, 22).Value = 4
$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 12).Value = "stimuli/img_h1yyu.png"
$ws.Cells.Item(32, 13).Value = 64.85294117647059
$ws.Cells.Item(32, 14).Value = 46.61764705882353
$ws.Cells.Item(32, 15).Value = 55.73529411764706
$ws.Cells.Item(32, 16).Value = 34
$ws.Cells.Item(32, 17).Value = 4
$ws.Cells.Item(32, 18).Value = 4
$ws.Cells.Item(32, 19).Value = 4
$ws.Cells.Item(32, 20).Value = 4
$ws.Cells.Item(32, 21).Value = 4
$ws.Cells.Item(32, 22).Value = 4
$ws.Cells.Item(33, 3).Value = 2
$ws.Cells.Item(33, 8).Value = "living_rooms"
$ws.Cells.Item(33, 12).Value = "stimuli/img_fmgjx.png"
$ws.Cells.Item(33, 13).Value = 79.90000000000001
$ws.Cells.Item(33, 14).Value = 56.975
$ws.Cells.Item(33, 15).Value = 68.4375
$ws.Cells.Item(33, 16).Value = 40
$ws.Cells.Item(33, 17).Value = 7
$ws.Cells.Item(33, 18).Value = 7
$ws.Cells.Item(33, 19).Value = 7
$ws.Cells.Item(33, 20).Value = 7
$ws.Cells.Item(33, 21).Value = 7
$ws.Cells.Item(33, 22).Value = 7
$ws.Cells.Item(34, 3).Value = 2
$ws.Cells.Item(34, 8).Value = "kitchens"
$ws.Cells.Item(34, 9).Value = "target"
$ws.Cells.Item(34, 11).Value = "j"
$ws.Cells.Item(34, 12).Value = "stimuli/img_oz18d.png"
$ws.Cells.Item(34, 13).Value = 78.93939393939394
$ws.Cells.Item(34, 14).Value = 61.03030303030303
$ws.Cells.Item(34, 15).Value = 69.98484848484848
$ws.Cells.Item(34, 16).Value = 33
$ws.Cells.Item(34, 17).Value = 8
$ws.Cells.Item(34, 18).Value = 8
$ws.Cells.Item(34, 19).Value = 8
$ws.Cells.Item(34, 20).Value = 8
$ws.Cells.Item(34, 21).Value = 8
$ws.Cells.Item(34, 22).Value = 8
$ws.Cells.Item(35, 3).Value = 2
$ws.Cells.Item(35, 8).Value = "living_rooms"
$ws.Cells.Item(35, 9).Value = "distractor"
$ws.Cells.Item(35, 11).Value = "f"
$ws.Cells.Item(35, 12).Value = "stimuli/img_abobq.png"
$ws.Cells.Item(35, 13).Value = 75.18421052631579
$ws.Cells.Item(35, 14).Value = 54.13157894736842
$ws.Cells.Item(35, 15).Value = 64.65789473684211
$ws.Cells.Item(35, 16).Value = 38
$ws.Cells.Item(35, 17).Value = 6
$ws.Cells.Item(35, 18).Value = 6
$ws.Cells.Item(35, 19).Value = 6
$ws.Cells.Item(35, 20).Value = 6
$ws.Cells.Item(35, 21).Value = 6
$ws.Cells.Item(35, 22).Value = 6
$ws.Cells.Item(36, 3).Value = 2
$ws.Cells.Item(36, 12).Value = "stimuli/img_t1cr9.png"
$ws.Cells.Item(36, 13).Value = 73.66666666666667
$ws.Cells.Item(36, 14).Value = 53.51515151515152
$ws.Cells.Item(36, 15).Value = 63.59090909090909
$ws.Cells.Item(36, 16).Value = 33
$ws.Cells.Item(36, 17).Value = 6
$ws.Cells.Item(36, 18).Value = 6
$ws.Cells.Item(36, 19).Value = 6
$ws.Cells.Item(36, 20).Value = 6
$ws.Cells.Item(36, 21).Value = 6
$ws.Cells.Item(36, 22).Value = 6
$ws.Cells.Item(37, 3).Value = 2
$ws.Cells.Item(37, 12).Value = "stimuli/img_57os5.png"
$ws.Cells.Item(37, 13).Value = 82.70588235294117
$ws.Cells.Item(37, 14).Value = 65.73529411764706
$ws.Cells.Item(37, 15).Value = 74.22058823529412
$ws.Cells.Item(37, 16).Value = 34
$ws.Cells.Item(37, 17).Value = 9
$ws.Cells.Item(37, 18).Value = 9
$ws.Cells.Item(37, 19).Value = 9
$ws.Cells.Item(37, 20).Value = 9
$ws.Cells.Item(37, 21).Value = 9
$ws.Cells.Item(37, 22).Value = 9
$ws.Cells.Item(38, 3).Value = 2
$ws.Cells.Item(38, 12).Value = "stimuli/img_lpr0l.png"
$ws.Cells.Item(38, 13).Value = 77.04651162790698
$ws.Cells.Item(38, 14).Value = 59.86046511627907
$ws.Cells.Item(38, 15).Value = 68.45348837209303
$ws.Cells.Item(38, 16).Value = 43
$ws.Cells.Item(38, 17).Value = 7
$ws.Cells.Item(38, 18).Value = 7
$ws.Cells.Item(38, 19).Value = 7
$ws.Cells.Item(38, 20).Value = 7
$ws.Cells.Item(38, 21).Value = 7
$ws.Cells.Item(38, 22).Value = 7
$ws.Cells.Item(39, 3).Value = 2
$ws.Cells.Item(39, 8).Value = "living_rooms"
$ws.Cells.Item(39, 9).Value = "distractor"
$ws.Cells.Item(39, 11).Value = "f"
$ws.Cells.Item(39, 12).Value = "stimuli/img_qrc78.png"
$ws.Cells.Item(39, 13).Value = 76.2
$ws.Cells.Item(39, 14).Value = 59.875
$ws.Cells.Item(39, 15).Value = 68.03749999999999
$ws.Cells.Item(39, 16).Value = 40
$ws.Cells.Item(39, 17).Value = 7
$ws.Cells.Item(39, 18).Value = 7
$ws.Cells.Item(39, 19).Value = 7
$ws.Cells.Item(39, 20).Value = 7
$ws.Cells.Item(39, 21).Value = 6
$ws.Cells.Item(39, 22).Value = 7
$ws.Cells.Item(40, 3).Value = 2
$ws.Cells.Item(40, 8).Value = "kitchens"
$ws.Cells.Item(40, 9).Value = "target"
$ws.Cells.Item(40, 11).Value = "j"
$ws.Cells.Item(40, 12).Value = "stimuli/img_z5osu.png"
$ws.Cells.Item(40, 13).Value = 71.42857142857143
$ws.Cells.Item(40, 14).Value = 47.34285714285714
$ws.Cells.Item(40, 15).Value = 59.38571428571429
$ws.Cells.Item(40, 16).Value = 35
$ws.Cells.Item(40, 17).Value = 5
$ws.Cells.Item(40, 18).Value = 5
$ws.Cells.Item(40, 19).Value = 5
$ws.Cells.Item(40, 20).Value = 5
$ws.Cells.Item(40, 21).Value = 5
$ws.Cells.Item(40, 22).Value = 5
$ws.Cells.Item(41, 3).Value = 2
$ws.Cells.Item(41, 12).Value = "stimuli/img_zh8ms.png"
$ws.Cells.Item(41, 13).Value = 59.82608695652174
$ws.Cells.Item(41, 14).Value = 39.43478260869565
$ws.Cells.Item(41, 15).Value = 49.6304347826087
$ws.Cells.Item(41, 16).Value = 46
$ws.Cells.Item(41, 17).Value = 4
$ws.Cells.Item(41, 18).Value = 4
$ws.Cells.Item(41, 19).Value = 4
$ws.Cells.Item(41, 20).Value = 4
$ws.Cells.Item(41, 21).Value = 4
$ws.Cells.Item(41, 22).Value = 4

Write-Host "Updated 537 cells"
